$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 138 (existing rows 138:149 shift down to 141:152)
$ws.Rows("138:140").Insert()

# --- New row 138 (Extra) ---
$ws.Range("A138").Value = 7
$ws.Range("B138").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C138").Value = "Ñuble"
$ws.Range("D138").Value = 44585
$ws.Range("E138").Value = 16
$ws.Range("F138").Value = 100112028
$ws.Range("G138").Value = "Sandia"
$ws.Range("H138").Value = "Sin especificar"
$ws.Range("I138").Value = "Extra"
$ws.Range("J138").Value = 200
$ws.Range("K138").Value = 2500
$ws.Range("L138").Value = 2500
$ws.Range("M138").Value = 2500
$ws.Range("N138").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O138").Value = "Región de O'Higgins"
$ws.Range("P138").Value = 2500
$ws.Range("Q138").Value = 1
$ws.Range("R138").Value = "Hortaliza"

# --- New row 139 (Primera) ---
$ws.Range("A139").Value = 7
$ws.Range("B139").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C139").Value = "Ñuble"
$ws.Range("D139").Value = 44585
$ws.Range("E139").Value = 16
$ws.Range("F139").Value = 100112028
$ws.Range("G139").Value = "Sandia"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 300
$ws.Range("K139").Value = 2000
$ws.Range("L139").Value = 2300
$ws.Range("M139").Value = 2150
$ws.Range("N139").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O139").Value = "Región de O'Higgins"
$ws.Range("P139").Value = 2150
$ws.Range("Q139").Value = 1
$ws.Range("R139").Value = "Hortaliza"

# --- New row 140 (Segunda) ---
$ws.Range("A140").Value = 7
$ws.Range("B140").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C140").Value = "Ñuble"
$ws.Range("D140").Value = 44585
$ws.Range("E140").Value = 16
$ws.Range("F140").Value = 100112028
$ws.Range("G140").Value = "Sandia"
$ws.Range("H140").Value = "Sin especificar"
$ws.Range("I140").Value = "Segunda"
$ws.Range("J140").Value = 200
$ws.Range("K140").Value = 1500
$ws.Range("L140").Value = 1800
$ws.Range("M140").Value = 1650
$ws.Range("N140").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O140").Value = "Región de O'Higgins"
$ws.Range("P140").Value = 1650
$ws.Range("Q140").Value = 1
$ws.Range("R140").Value = "Hortaliza"
